$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 7.805142343733317
$arr[1,0] = 7.734595247841789
$arr[2,0] = 7.692877315427562
$arr[3,0] = 7.676296997732809
$arr[4,0] = 7.673569735895787
$arr[5,0] = 7.692651983142277
$arr[6,0] = 7.780497128041432
$arr[7,0] = 7.964608397873418
$arr[8,0] = 8.105938624933682
$arr[9,0] = 8.171295309601078
$arr[10,0] = 8.196175663739321
$arr[11,0] = 8.19081174847636
$arr[12,0] = 8.173339727560677
$arr[13,0] = 8.16265404489536
$arr[14,0] = 8.101686965581935
$arr[15,0] = 8.064542783282793
$arr[16,0] = 8.043279979764334
$arr[17,0] = 8.036098868334573
$arr[18,0] = 8.068486493878371
$arr[19,0] = 8.178468297784235
$arr[20,0] = 8.251102328675323
$arr[21,0] = 8.212274376025592
$arr[22,0] = 8.066703254464299
$arr[23,0] = 7.913650899888146
$ws.Range("B2:B25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 6.31665926788368
$arr[1,0] = 6.337819254527891
$arr[2,0] = 6.351854287584858
$arr[3,0] = 6.357834381788795
$arr[4,0] = 6.358843074774499
$arr[5,0] = 6.351933883581694
$arr[6,0] = 6.323737904895206
$arr[7,0] = 6.276787924866975
$arr[8,0] = 6.247476662488154
$arr[9,0] = 6.235289585132072
$arr[10,0] = 6.230841469748504
$arr[11,0] = 6.231791997636462
$arr[12,0] = 6.234920280301628
$arr[13,0] = 6.236858233948134
$arr[14,0] = 6.248296358917672
$arr[15,0] = 6.25560836962906
$arr[16,0] = 6.259921862543867
$arr[17,0] = 6.261400796035262
$arr[18,0] = 6.25481882085319
$arr[19,0] = 6.233996883201201
$arr[20,0] = 6.221362121310024
$arr[21,0] = 6.228015784857124
$arr[22,0] = 6.255175434247369
$arr[23,0] = 6.288586199824935
$ws.Range("D2:D25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.35382138172144
$arr[1,0] = 15.42957719045927
$arr[2,0] = 14.83910729246217
$arr[3,0] = 14.59298196726907
$arr[4,0] = 14.55178968477319
$arr[5,0] = 14.83580984431916
$arr[6,0] = 16.04006048101398
$arr[7,0] = 18.25237702746207
$arr[8,0] = 19.89630835299693
$arr[9,0] = 20.60323771848714
$arr[10,0] = 20.8650934224524
$arr[11,0] = 20.8089573694548
$arr[12,0] = 20.62489762776404
$arr[13,0] = 20.51139592484072
$arr[14,0] = 19.84928909168439
$arr[15,0] = 19.43265928316738
$arr[16,0] = 19.1891725196347
$arr[17,0] = 19.10606926505558
$arr[18,0] = 19.47740865672518
$arr[19,0] = 20.67911870755498
$arr[20,0] = 21.43046552539449
$arr[21,0] = 21.03255949208876
$arr[22,0] = 19.45718980468327
$arr[23,0] = 17.64443254072353
$ws.Range("E2:E25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 36.05238240123538
$arr[1,0] = 35.59185804786794
$arr[2,0] = 35.3182212984111
$arr[3,0] = 35.20912487836423
$arr[4,0] = 35.19115854513333
$arr[5,0] = 35.31674006564246
$arr[6,0] = 35.89177168751594
$arr[7,0] = 37.0863905304528
$arr[8,0] = 37.99735719480683
$arr[9,0] = 38.41735407425632
$arr[10,0] = 38.57706132556309
$arr[11,0] = 38.54263822809935
$arr[12,0] = 38.43048078952397
$arr[13,0] = 38.36186337513907
$arr[14,0] = 37.97000990799291
$arr[15,0] = 37.73095267456655
$arr[16,0] = 37.59398816918288
$arr[17,0] = 37.54771051536291
$arr[18,0] = 37.75634638928756
$arr[19,0] = 38.46340725536636
$arr[20,0] = 38.9293158663849
$arr[21,0] = 38.68034979881199
$arr[22,0] = 37.74486440786814
$arr[23,0] = 36.7568122842541
$ws.Range("F2:F25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 3.660618554858681
$arr[1,0] = 3.665194562076497
$arr[2,0] = 3.668142403746693
$arr[3,0] = 3.669378578032615
$arr[4,0] = 3.669585956819262
$arr[5,0] = 3.668158933677465
$arr[6,0] = 3.662167787640302
$arr[7,0] = 3.651507824465597
$arr[8,0] = 3.644329032281586
$arr[9,0] = 3.641202741195801
$arr[10,0] = 3.640038759955966
$arr[11,0] = 3.640288562722607
$arr[12,0] = 3.641106582334376
$arr[13,0] = 3.641610226677213
$arr[14,0] = 3.644536133296092
$arr[15,0] = 3.646366662925568
$arr[16,0] = 3.64743266391704
$arr[17,0] = 3.64779585342994
$arr[18,0] = 3.646170442337988
$arr[19,0] = 3.6408657719968
$arr[20,0] = 3.637514646711323
$arr[21,0] = 3.639292665791555
$arr[22,0] = 3.646259111276472
$arr[23,0] = 3.654276172432507
$ws.Range("G2:G25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.48867995424968
$arr[1,0] = 15.9711343104385
$arr[2,0] = 15.65076542929826
$arr[3,0] = 15.51978843109278
$arr[4,0] = 15.4980211179643
$arr[5,0] = 15.64900043331175
$arr[6,0] = 16.31090111012266
$arr[7,0] = 17.57902002748707
$arr[8,0] = 18.48118027224943
$arr[9,0] = 18.88312045810357
$arr[10,0] = 19.03395879633925
$arr[11,0] = 19.00153613463999
$arr[12,0] = 18.89555824127583
$arr[13,0] = 18.83046138099199
$arr[14,0] = 18.45472903696972
$arr[15,0] = 18.22195042842245
$arr[16,0] = 18.08727433993865
$arr[17,0] = 18.04154487315979
$arr[18,0] = 18.24681281124231
$arr[19,0] = 18.926724765462
$arr[20,0] = 19.36305534339105
$arr[21,0] = 19.13095836584036
$arr[22,0] = 18.23557515938923
$arr[23,0] = 17.24038970069342
$ws.Range("K2:K25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 19.81168637808847
$arr[1,0] = 19.86211620128597
$arr[2,0] = 19.89500083856163
$arr[3,0] = 19.90888395235965
$arr[4,0] = 19.91121834974602
$arr[5,0] = 19.89518611875938
$arr[6,0] = 19.8286756264874
$arr[7,0] = 19.71351308494455
$arr[8,0] = 19.63824460439897
$arr[9,0] = 19.60603993653978
$arr[10,0] = 19.59413839785163
$arr[11,0] = 19.5966885297767
$arr[12,0] = 19.60505489837225
$arr[13,0] = 19.61021781539949
$arr[14,0] = 19.64039029504325
$arr[15,0] = 19.65942212357049
$arr[16,0] = 19.6705602419125
$arr[17,0] = 19.67436428061513
$arr[18,0] = 19.65737632495387
$arr[19,0] = 19.60258951485256
$arr[20,0] = 19.56849539460517
$arr[21,0] = 19.58653505968358
$arr[22,0] = 19.65830061863337
$arr[23,0] = 19.74302990214388
$ws.Range("N2:N25").Value = $arr

Write-Host "applied loading_percent updates"